$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A71").Value = "Create a country"
$ws.Range("B71").Value = "PASSED"
$ws.Range("C71").Value = "chrome"
$ws.Range("D71").NumberFormat = "@"
$ws.Range("D71").Value = "12.10.22"
$ws.Range("D71").Style = "Normal"

$ws.Range("A72").Value = "Create a country  with parameter data"
$ws.Range("B72").Value = "PASSED"
$ws.Range("C72").Value = "chrome"
$ws.Range("D72").NumberFormat = "@"
$ws.Range("D72").Value = "12.10.22"
$ws.Range("D72").Style = "Normal"
